$wb = $excel.ActiveWorkbook

# --- Clients sheet: fill "Test Passed" column (E) with TRUE for all data rows ---
$wsClients = $wb.Worksheets.Item("Clients")
$wsClients.Range("E2:E51").Value = $true
$wsClients.Range("E52").Select()

# --- Projects sheet: fill "Test Passed" column (G) with TRUE for all data rows ---
# Do this last so Projects ends up as the active sheet/tab (matches final state).
$wsProjects = $wb.Worksheets.Item("Projects")
$wsProjects.Activate()
$wsProjects.Range("G2:G46").Value = $true
$wsProjects.Range("G2:G46").Select()
